$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the placeholder "router/switch/server" rows (4-8, previously mostly empty)
# with the real purchased items, and fill in the remaining rows (9-12) that were
# blank before. Quantities/unit prices are entered so that the existing shared
# formula D4:D12 (=B*C) and the running total I2 (=SUM(D2:D30)) recompute
# automatically.
$ws.Range("A4").Value = "cisco isr 4331"
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 61965

$ws.Range("A5").Value = "cisco ws-c2960x 24psq-l"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 366363

$ws.Range("A6").Value = "cisco c9300-48ub-a"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 5544572

$ws.Range("A7").Value = "dahua lm22-b200s"
$ws.Range("B7").Value = 7
$ws.Range("C7").Value = 28500

$ws.Range("A8").Value = "dell optiplex 3020 sff"
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 29900

$ws.Range("A9").Value = "cisco cbw140ac"
$ws.Range("B9").Value = 23
$ws.Range("C9").Value = 41090

$ws.Range("A10").Value = "white shark cryus"
$ws.Range("B10").Value = 7
$ws.Range("C10").Value = 4999

$ws.Range("A11").Value = "white shark commandos tkl"
$ws.Range("B11").Value = 7
$ws.Range("C11").Value = 15490

$ws.Range("B12").Value = 1

# Widen column A so the longer product names fit (target stored width 35.109375
# chars; the COM layer quantizes ColumnWidth to 1/6-character steps, so
# 34.3333... is the closest input that rounds to the nearest achievable width).
$ws.Columns.Item(1).ColumnWidth = 34.333333333333336

# Match the saved cursor/selection position.
$ws.Range("F4").Select() | Out-Null
